# Error.xlsx - "add the scene chess move effect"
# Error code 2000 used to mean "需要等级到达{0}级" (level requirement not met).
# It is repurposed to report a new failure condition for entering a map/scene
# via the chess-move effect: "不符合地图进入条件" (map entry condition not met).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Error")

# Row 17 -> A17=2000, B17 holds the description shown to players for that code.
$ws.Range("B17").Value = "不符合地图进入条件"

# Leave the cursor on the cell that was just edited, matching the authored
# workbook's last-saved selection.
[void]$ws.Range("B17").Select()

Write-Host "Updated Error code 2000 (B17) to the new map-entry error message."
